# final changes to questionnaires and video script
#
# Slide 11 ("clean the gutters") - Content Placeholder 2:
#   - paragraph 1: split " remembers that you cleaned the gutters last
#     winter, " into " remembers that " / "last winter y" /
#     "ou completed the task of cleaning " / "the " / "gutters, "
#   - paragraph 3: merge "winter " + "with " + "one click" -> "winter with one click"
#   - paragraph 6: merge "find " + "out which ..." -> "find out which ..."
#   - paragraph 7: merge "share " + "your experience ..." -> "share your experience ..."
#
# Slide 13 ("learns about recurring patterns") - Content Placeholder 2:
#   - paragraph 9: merge "learns about recurring patterns in your life and "
#     + "what " + "you're likely to do in the future" into a single run

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)
$tr11 = $sh11.TextFrame.TextRange

# Paragraph 1: "zaplify remembers that you cleaned the gutters last winter, so it can..."
$para1 = $tr11.Paragraphs(1, 1)
$replaceStart = $para1.Start + 7
$replaceLen = (" remembers that you cleaned the gutters last winter, ").Length
$oldRun = $tr11.Characters($replaceStart, $replaceLen)
$oldRun.Text = " remembers that last winter you completed the task of cleaning the gutters, "

$pos = $para1.Start + 7
$seg1 = $tr11.Characters($pos, 16)
$seg1.Text = $seg1.Text
$pos = $pos + 16
$seg2 = $tr11.Characters($pos, 13)
$seg2.Text = $seg2.Text
$pos = $pos + 13
$seg3 = $tr11.Characters($pos, 34)
$seg3.Text = $seg3.Text
$pos = $pos + 34
$seg4 = $tr11.Characters($pos, 4)
$seg4.Text = $seg4.Text
$pos = $pos + 4
$seg5 = $tr11.Characters($pos, 9)
$seg5.Text = $seg5.Text

# Paragraph 3: "allow you to call back the service you used last winter with one click"
$para3 = $tr11.Paragraphs(3, 1)
$mergeStart = $para3.Start + 49
$mergeLen = 21
$mergeRange = $tr11.Characters($mergeStart, $mergeLen)
$mergeRange.Text = "winter with one click"

# Paragraph 6: "find out which services your friends have used and how they rated them"
$para6 = $tr11.Paragraphs(6, 1)
$full6 = $tr11.Characters($para6.Start, $para6.Length)
$full6.Text = "find out which services your friends have used and how they rated them"

# Paragraph 7: "share your experience with your friends"
$para7 = $tr11.Paragraphs(7, 1)
$full7 = $tr11.Characters($para7.Start, $para7.Length)
$full7.Text = "share your experience with your friends"

# ---------------------------------------------------------------------------
# Slide 13
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange

# Paragraph 9: "zaplify learns about recurring patterns in your life and what you're likely to do in the future"
$para9 = $tr13.Paragraphs(9, 1)
$mergeStart9 = $para9.Start + 8
$mergeLen9 = 87
$mergeRange9 = $tr13.Characters($mergeStart9, $mergeLen9)
$mergeRange9.Text = "learns about recurring patterns in your life and what you’re likely to do in the future"
